$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format before writing, so numeric-looking
# strings (e.g. "1.008") are stored as text rather than being parsed as numbers,
# matching the original inlineStr text cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.034.86'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '1.830.23'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '311.61'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D7').Value = '0.4648'
$ws.Range('E7').Value = '  -1.93%  '
$ws.Range('D8').Value = '0.3710'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = '0.07406'
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').Value = '0.8662'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').Value = '20.01'
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('D12').Value = '0.07855'
$ws.Range('E12').Value = '  +7.49%  '
$ws.Range('D13').Value = '1.835.05'
$ws.Range('E13').Value = '  -4.59%  '
$ws.Range('D14').Value = '6.637'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '5.362'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '92.12'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '0.000009073'
$ws.Range('E18').Value = '  +3.13%  '
$ws.Range('D20').Value = '14.68'
$ws.Range('D21').Value = '27.073.90'
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = '5.167'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '2.065.57'
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('D25').Value = '152.78'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').Value = '1.841'
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('D27').Value = '18.27'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('D28').Value = '2.101'
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('D29').Value = '5.137'
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('D30').Value = '115.78'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').Value = '0.08873'
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').Value = '2.985'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('D33').Value = '0.7296'
$ws.Range('E33').Value = '  -2.54%  '
$ws.Range('D34').Value = '4.452'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').Value = '1.134'
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('D36').Value = '2.479'
$ws.Range('E36').Value = '  +2.61%  '
$ws.Range('D37').Value = '1.082'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = '0.01952'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').Value = '7.380'
$ws.Range('E39').Value = '  +2.14%  '
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').Value = '2.928'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').Value = '0.5178'
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('D43').Value = '0.1632'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').Value = '0.8576'
$ws.Range('E44').Value = '  -15.11%  '
$ws.Range('D45').Value = '8.238'
$ws.Range('E45').Value = '  -2.90%  '
$ws.Range('D46').Value = '0.4843'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.008'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.24'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('D49').Value = '102.79'
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('D50').Value = '1.624'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('D51').Value = '0.06251'
$ws.Range('E51').Value = '  -0.82%  '

# Restore default (Normal) style on the touched columns so no stray
# cell-format/style index is left behind by the NumberFormat tweak above.
$ws.Range("D2:E51").Style = "Normal"

